$wb = $excel.ActiveWorkbook

# Sheet 1 (展览) - column F ("想去人数") value updates
$ws1 = $wb.Worksheets.Item(1)
$ws1.Cells.Item(2, 6).Value = 807
$ws1.Cells.Item(3, 6).Value = 549
$ws1.Cells.Item(5, 6).Value = 507
$ws1.Cells.Item(8, 6).Value = 47
$ws1.Cells.Item(9, 6).Value = 122
$ws1.Cells.Item(11, 6).Value = 1185
$ws1.Cells.Item(14, 6).Value = 847
$ws1.Cells.Item(15, 6).Value = 853
$ws1.Cells.Item(20, 6).Value = 739
$ws1.Cells.Item(21, 6).Value = 1730
$ws1.Cells.Item(22, 6).Value = 2656
$ws1.Cells.Item(23, 6).Value = 753
$ws1.Cells.Item(24, 6).Value = 79
$ws1.Cells.Item(25, 6).Value = 2024
$ws1.Cells.Item(26, 6).Value = 497
$ws1.Cells.Item(27, 6).Value = 2922
$ws1.Cells.Item(28, 6).Value = 546
$ws1.Cells.Item(30, 6).Value = 86
$ws1.Cells.Item(31, 6).Value = 712
$ws1.Cells.Item(33, 6).Value = 119
$ws1.Cells.Item(35, 6).Value = 1020
$ws1.Cells.Item(36, 6).Value = 1739
$ws1.Cells.Item(37, 6).Value = 366
$ws1.Cells.Item(39, 6).Value = 546
$ws1.Cells.Item(40, 6).Value = 174
$ws1.Cells.Item(41, 6).Value = 129
$ws1.Cells.Item(42, 6).Value = 163
$ws1.Cells.Item(43, 6).Value = 33

# Sheet 2 (演出) - column F updates
$ws2 = $wb.Worksheets.Item(2)
$ws2.Cells.Item(9, 6).Value = 11

# Sheet 4 (全部类型) - column F updates
$ws4 = $wb.Worksheets.Item(4)
$ws4.Cells.Item(3, 6).Value = 807
$ws4.Cells.Item(4, 6).Value = 549
$ws4.Cells.Item(6, 6).Value = 507
$ws4.Cells.Item(9, 6).Value = 47
$ws4.Cells.Item(10, 6).Value = 122
$ws4.Cells.Item(12, 6).Value = 1185
$ws4.Cells.Item(14, 6).Value = 847
$ws4.Cells.Item(15, 6).Value = 853
$ws4.Cells.Item(21, 6).Value = 739
$ws4.Cells.Item(22, 6).Value = 1730
$ws4.Cells.Item(23, 6).Value = 2656
$ws4.Cells.Item(24, 6).Value = 753
$ws4.Cells.Item(25, 6).Value = 79
$ws4.Cells.Item(28, 6).Value = 2922
$ws4.Cells.Item(29, 6).Value = 546
$ws4.Cells.Item(31, 6).Value = 11
$ws4.Cells.Item(34, 6).Value = 86
$ws4.Cells.Item(36, 6).Value = 712
$ws4.Cells.Item(38, 6).Value = 119
$ws4.Cells.Item(40, 6).Value = 1020
$ws4.Cells.Item(41, 6).Value = 1739
$ws4.Cells.Item(43, 6).Value = 366
$ws4.Cells.Item(44, 6).Value = 546
$ws4.Cells.Item(45, 6).Value = 174
$ws4.Cells.Item(46, 6).Value = 129
$ws4.Cells.Item(47, 6).Value = 163
$ws4.Cells.Item(48, 6).Value = 33

